# feat: add 2022-Q1 data
#
# The former "总计" sheet (sheetId=5) is turned into the new "2022-Q1" fund
# holdings sheet, and a duplicate of it (sheetId=6) becomes the new "总计"
# summary sheet with the 2022-Q1 row folded in.

$wb = $excel.ActiveWorkbook

$templateFundSheet = $wb.Worksheets.Item(4)   # "2021-Q4" - fund-holdings layout template (B..H headers)
$oldTotalSheet      = $wb.Worksheets.Item(5)   # current "总计" sheet -> will become "2022-Q1"

# Duplicate the current "总计" sheet so the copy (placed right after it) can
# become the new "总计" sheet; this preserves sheetPr/sheetFormatPr/pageMargins
# and keeps the original sheet's identity (sheetId/relationship) intact for
# re-use as "2022-Q1".
$oldTotalSheet.Copy([System.Reflection.Missing]::Value, $oldTotalSheet)
$newTotalSheet = $wb.Worksheets.Item(6)

# ---------------------------------------------------------------------------
# Build the new "总计" sheet: old quarterly rows plus the new 2022-Q1 row
# ---------------------------------------------------------------------------
$newTotalSheet.Cells.Clear()
$oldTotalSheet.Range("B1:D1").Copy()
$newTotalSheet.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$oldTotalSheet.Range("A2").Copy()
$newTotalSheet.Range("A2:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newTotalSheet.Range("B1").Value = "日期"
$newTotalSheet.Range("C1").Value = "持有数量(只)"
$newTotalSheet.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 4, 0.5600000000000001),
    @(1, "2021-Q4", 6, 1.83),
    @(2, "2021-Q2", 4, 3.01),
    @(3, "2021-Q1", 2, 2.27),
    @(4, "2020-Q4", 3, 1.08)
)

$r = 2
foreach ($row in $totalRows) {
    $newTotalSheet.Cells.Item($r, 1).Value = $row[0]
    $newTotalSheet.Cells.Item($r, 2).Value = $row[1]
    $newTotalSheet.Cells.Item($r, 3).Value = $row[2]
    $newTotalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Turn the original "总计" sheet into the "2022-Q1" fund-holdings sheet
# ---------------------------------------------------------------------------
$oldTotalSheet.Cells.Clear()

# Copy header style (B1:H1) and index column style (A2:A5) from the 2021-Q4 template
$templateFundSheet.Range("B1:H1").Copy()
$oldTotalSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$templateFundSheet.Range("A2").Copy()
$oldTotalSheet.Range("A2:A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$oldTotalSheet.Range("B1").Value = "基金代码"
$oldTotalSheet.Range("C1").Value = "基金名称"
$oldTotalSheet.Range("D1").Value = "基金规模"
$oldTotalSheet.Range("E1").Value = "股票总仓位"
$oldTotalSheet.Range("F1").Value = "仓位占比"
$oldTotalSheet.Range("G1").Value = "持有市值(亿元)"
$oldTotalSheet.Range("H1").Value = "仓位排名"

# Use a scratch cell (far outside the sheet's used range) to stamp a "text"
# number format once, then copy its VALUE only onto the real destination
# cells, so the destination stays on the default style (no s= attribute)
# while the numeric-looking text is still stored as a text string - matching
# the source data's typing.
$scratch = $oldTotalSheet.Cells.Item(200, 200)

function Set-TextValue($cell, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}

$fundRows = @(
    @(0, "012210", "申万菱信智能汽车股票型证券投资基金A", "4.76", "82.52", "3.69", "0.1756", 10),
    @(1, "162102", "金鹰中小盘精选混合",                 "4.60", "76.52", "3.71", "0.1707", 5),
    @(2, "210009", "金鹰核心资源混合",                   "3.86", "94.96", "4.26", "0.1644", 10),
    @(3, "012211", "申万菱信智能汽车股票型证券投资基金C", "1.40", "82.52", "3.69", "0.0517", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $oldTotalSheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 2) $row[1]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 3) $row[2]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 4) $row[3]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 5) $row[4]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 6) $row[5]
    Set-TextValue $oldTotalSheet.Cells.Item($r, 7) $row[6]
    $oldTotalSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$scratch.Clear()
$excel.CutCopyMode = 0

# Rename the old sheet away from "总计" first so the new sheet can take that name
$oldTotalSheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# Restore the original active sheet/tab selection
$wb.Worksheets.Item(1).Activate()
